$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with updated crypto market
# data. Values are entered with a leading apostrophe so Excel stores them as
# literal text (matching the existing text-formatted Price/Volume columns)
# rather than auto-converting them to numbers/percentages; the cell style is
# then reset to Normal so the quote-prefix indicator does not linger.

$ws.Range("D2").Value = "'313.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'9.09%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'10.30%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.336"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.34%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07659"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'14.36%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.891"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'7.50%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.749"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'10.14%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'19.05%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9194"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.94%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.01770"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2,637.52%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1725"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'8.82%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07585"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'12.56%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08319"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'7.85%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03039"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.77%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09906"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'10.33%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001525"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.44%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04569"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.88%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006215"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.77%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.468"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.74%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.244"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.10%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'3.44%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1334"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.95%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.240"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'4.24%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'0.001221"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.55%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'9.53%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001298"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'8.27%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'7.47%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04648"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'9.08%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007193"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'7.26%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1373"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'10.88%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002257"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.14%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01443"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'8.66%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006198"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'9.22%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-3.82%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.65%"
$ws.Range("E47").Style = "Normal"
